# Insert a new data row at row 12 (pushing the existing rows 12-126 down to 13-127)
# and populate the new row 12 with the latest weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 12, shifting rows 12-126 down to 13-127.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new record's data.
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(12, 3).Value = "Bíobío"
$ws.Cells.Item(12, 4).Value = 44685
$ws.Cells.Item(12, 5).Value = 8
$ws.Cells.Item(12, 6).Value = 100112043
$ws.Cells.Item(12, 7).Value = "Pepino ensalada"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 220
$ws.Cells.Item(12, 11).Value = 17000
$ws.Cells.Item(12, 12).Value = 18000
$ws.Cells.Item(12, 13).Value = 17455
$ws.Cells.Item(12, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 291
$ws.Cells.Item(12, 17).Value = 60
$ws.Cells.Item(12, 18).Value = "Hortaliza"
